$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 16.28844733333333
$ws.Range("N2").Value = 48.865342
$ws.Range("O2").Value = 0.2176904746803693
$ws.Range("P2").Value = 0.2176904746803693
$ws.Range("Q2").Value = 2.130604923954222
$ws.Range("R2").Value = 19.175444315588
$ws.Range("S2").Value = 0.2176904746803693
$ws.Range("T2").Value = 0.2176904746803693

# Row 3
$ws.Range("M3").Value = 27.61090666666666
$ws.Range("N3").Value = 82.83271999999999
$ws.Range("O3").Value = 0.3690119294748028
$ws.Range("P3").Value = 0.3690119294748029
$ws.Range("Q3").Value = 3.611635442897777
$ws.Range("R3").Value = 32.50471898607999
$ws.Range("S3").Value = 0.3690119294748028
$ws.Range("T3").Value = 0.3690119294748029

# Row 4
$ws.Range("M4").Value = 26.266325
$ws.Range("N4").Value = 78.798975
$ws.Range("O4").Value = 0.3510419771967738
$ws.Range("P4").Value = 0.3510419771967739
$ws.Range("Q4").Value = 3.435757886183333
$ws.Range("R4").Value = 30.92182097565
$ws.Range("S4").Value = 0.3510419771967738
$ws.Range("T4").Value = 0.3510419771967739

# Row 5
$ws.Range("M5").Value = 4.658207333333333
$ws.Range("N5").Value = 13.974622
$ws.Range("O5").Value = 0.06225561864805391
$ws.Range("P5").Value = 0.06225561864805392
$ws.Range("Q5").Value = 0.6093152575008888
$ws.Range("R5").Value = 5.483837317508
$ws.Range("S5").Value = 0.06225561864805391
$ws.Range("T5").Value = 0.06225561864805392
